# Reorder devices in the benchmark table: swap the "Alveo U200" block
# (rows 2-5) with the "Alveo U280" block (rows 10-13). Only the raw data
# columns are swapped (A, B, E, F, H, I, J, L); the formula columns
# (C, D, G, K) stay bound to their own row and simply recalculate.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @("A","B","E","F","H","I","J","L")
$rowPairs = @(@(2,10), @(3,11), @(4,12), @(5,13))

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    foreach ($col in $dataCols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

$excel.Calculate()

$ws.Range("E22").Select() | Out-Null
